# Auto-generated: relocates the full records (columns A:AY) that lived in
# rows 9-15 to their new row positions, per the mapping taken from the diff.
# Every destination cell is written explicitly with its exact original type
# (number / boolean / text) so that values which merely look numeric or look
# like a date (e.g. the text "1" or the text "2022-08-22") are not silently
# reinterpreted by Excel's COM auto-detection when assigned through .Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- row 9 ----
$ws.Cells.Item(9, 1).Value = 103636871  # A9
$ws.Cells.Item(9, 2).Value = 56395  # B9
$ws.Cells.Item(9, 3).Formula = '''Ovaliderad'  # C9
$ws.Cells.Item(9, 4).Formula = '''NT'  # D9
$ws.Cells.Item(9, 5).Value = 100109  # E9
$ws.Cells.Item(9, 6).Formula = '''Tretåig hackspett'  # F9
$ws.Cells.Item(9, 7).Formula = '''Picoides tridactylus'  # G9
$ws.Cells.Item(9, 8).Formula = '''(Linnaeus, 1758)'  # H9
$ws.Cells.Item(9, 9).Formula = ''''  # I9
$ws.Cells.Item(9, 10).Value = $null  # J9
$ws.Cells.Item(9, 11).Formula = ''''  # K9
$ws.Cells.Item(9, 12).Formula = ''''  # L9
$ws.Cells.Item(9, 13).Formula = '''äldre spår'  # M9
$ws.Cells.Item(9, 14).Formula = ''''  # N9
$ws.Cells.Item(9, 15).Value = $null  # O9
$ws.Cells.Item(9, 16).Formula = '''Bergom - Rödön, Jmt'  # P9
$ws.Cells.Item(9, 17).Value = 472516.7485192241  # Q9
$ws.Cells.Item(9, 18).Value = 7016947.556961586  # R9
$ws.Cells.Item(9, 19).Value = 10  # S9
$ws.Cells.Item(9, 20).Formula = '''Jämtland'  # T9
$ws.Cells.Item(9, 21).Formula = '''Krokom'  # U9
$ws.Cells.Item(9, 22).Formula = '''Jämtland'  # V9
$ws.Cells.Item(9, 23).Formula = '''Rödön'  # W9
$ws.Cells.Item(9, 24).Value = $null  # X9
$ws.Cells.Item(9, 25).Formula = '''2022-09-08'  # Y9
$ws.Cells.Item(9, 26).Formula = '''00:00'  # Z9
$ws.Cells.Item(9, 27).Formula = '''2022-09-08'  # AA9
$ws.Cells.Item(9, 28).Formula = '''00:00'  # AB9
$ws.Cells.Item(9, 29).Value = $null  # AC9
$ws.Cells.Item(9, 30).Value = $false  # AD9
$ws.Cells.Item(9, 31).Value = $false  # AE9
$ws.Cells.Item(9, 32).Value = $null  # AF9
$ws.Cells.Item(9, 33).Value = $false  # AG9
$ws.Cells.Item(9, 34).Value = $null  # AH9
$ws.Cells.Item(9, 35).Value = $null  # AI9
$ws.Cells.Item(9, 36).Value = $null  # AJ9
$ws.Cells.Item(9, 37).Value = $null  # AK9
$ws.Cells.Item(9, 38).Value = $null  # AL9
$ws.Cells.Item(9, 39).Value = $null  # AM9
$ws.Cells.Item(9, 40).Value = $null  # AN9
$ws.Cells.Item(9, 41).Value = $null  # AO9
$ws.Cells.Item(9, 42).Value = $null  # AP9
$ws.Cells.Item(9, 43).Value = $null  # AQ9
$ws.Cells.Item(9, 44).Value = $null  # AR9
$ws.Cells.Item(9, 45).Value = $null  # AS9
$ws.Cells.Item(9, 46).Formula = ''''  # AT9
$ws.Cells.Item(9, 47).Value = $null  # AU9
$ws.Cells.Item(9, 48).Value = $null  # AV9
$ws.Cells.Item(9, 49).Formula = '''Benny Öwre'  # AW9
$ws.Cells.Item(9, 50).Formula = '''Benny Öwre'  # AX9
$ws.Cells.Item(9, 51).Formula = ''''  # AY9

# ---- row 10 ----
$ws.Cells.Item(10, 1).Value = 103636870  # A10
$ws.Cells.Item(10, 2).Value = 56395  # B10
$ws.Cells.Item(10, 3).Formula = '''Ovaliderad'  # C10
$ws.Cells.Item(10, 4).Formula = '''NT'  # D10
$ws.Cells.Item(10, 5).Value = 100109  # E10
$ws.Cells.Item(10, 6).Formula = '''Tretåig hackspett'  # F10
$ws.Cells.Item(10, 7).Formula = '''Picoides tridactylus'  # G10
$ws.Cells.Item(10, 8).Formula = '''(Linnaeus, 1758)'  # H10
$ws.Cells.Item(10, 9).Formula = ''''  # I10
$ws.Cells.Item(10, 10).Value = $null  # J10
$ws.Cells.Item(10, 11).Formula = ''''  # K10
$ws.Cells.Item(10, 12).Formula = ''''  # L10
$ws.Cells.Item(10, 13).Formula = '''äldre spår'  # M10
$ws.Cells.Item(10, 14).Formula = ''''  # N10
$ws.Cells.Item(10, 15).Value = $null  # O10
$ws.Cells.Item(10, 16).Formula = '''Bergom - Rödön, Jmt'  # P10
$ws.Cells.Item(10, 17).Value = 472507.7934395059  # Q10
$ws.Cells.Item(10, 18).Value = 7016954.384313107  # R10
$ws.Cells.Item(10, 19).Value = 10  # S10
$ws.Cells.Item(10, 20).Formula = '''Jämtland'  # T10
$ws.Cells.Item(10, 21).Formula = '''Krokom'  # U10
$ws.Cells.Item(10, 22).Formula = '''Jämtland'  # V10
$ws.Cells.Item(10, 23).Formula = '''Rödön'  # W10
$ws.Cells.Item(10, 24).Value = $null  # X10
$ws.Cells.Item(10, 25).Formula = '''2022-09-08'  # Y10
$ws.Cells.Item(10, 26).Formula = '''00:00'  # Z10
$ws.Cells.Item(10, 27).Formula = '''2022-09-08'  # AA10
$ws.Cells.Item(10, 28).Formula = '''00:00'  # AB10
$ws.Cells.Item(10, 29).Value = $null  # AC10
$ws.Cells.Item(10, 30).Value = $false  # AD10
$ws.Cells.Item(10, 31).Value = $false  # AE10
$ws.Cells.Item(10, 32).Value = $null  # AF10
$ws.Cells.Item(10, 33).Value = $false  # AG10
$ws.Cells.Item(10, 34).Value = $null  # AH10
$ws.Cells.Item(10, 35).Value = $null  # AI10
$ws.Cells.Item(10, 36).Value = $null  # AJ10
$ws.Cells.Item(10, 37).Value = $null  # AK10
$ws.Cells.Item(10, 38).Value = $null  # AL10
$ws.Cells.Item(10, 39).Value = $null  # AM10
$ws.Cells.Item(10, 40).Value = $null  # AN10
$ws.Cells.Item(10, 41).Value = $null  # AO10
$ws.Cells.Item(10, 42).Value = $null  # AP10
$ws.Cells.Item(10, 43).Value = $null  # AQ10
$ws.Cells.Item(10, 44).Value = $null  # AR10
$ws.Cells.Item(10, 45).Value = $null  # AS10
$ws.Cells.Item(10, 46).Formula = ''''  # AT10
$ws.Cells.Item(10, 47).Value = $null  # AU10
$ws.Cells.Item(10, 48).Value = $null  # AV10
$ws.Cells.Item(10, 49).Formula = '''Benny Öwre'  # AW10
$ws.Cells.Item(10, 50).Formula = '''Benny Öwre'  # AX10
$ws.Cells.Item(10, 51).Formula = ''''  # AY10

# ---- row 11 ----
$ws.Cells.Item(11, 1).Value = 103636872  # A11
$ws.Cells.Item(11, 2).Value = 56395  # B11
$ws.Cells.Item(11, 3).Formula = '''Ovaliderad'  # C11
$ws.Cells.Item(11, 4).Formula = '''NT'  # D11
$ws.Cells.Item(11, 5).Value = 100109  # E11
$ws.Cells.Item(11, 6).Formula = '''Tretåig hackspett'  # F11
$ws.Cells.Item(11, 7).Formula = '''Picoides tridactylus'  # G11
$ws.Cells.Item(11, 8).Formula = '''(Linnaeus, 1758)'  # H11
$ws.Cells.Item(11, 9).Formula = ''''  # I11
$ws.Cells.Item(11, 10).Value = $null  # J11
$ws.Cells.Item(11, 11).Formula = ''''  # K11
$ws.Cells.Item(11, 12).Formula = ''''  # L11
$ws.Cells.Item(11, 13).Formula = '''äldre spår'  # M11
$ws.Cells.Item(11, 14).Formula = ''''  # N11
$ws.Cells.Item(11, 15).Value = $null  # O11
$ws.Cells.Item(11, 16).Formula = '''Bergom - Rödön, Jmt'  # P11
$ws.Cells.Item(11, 17).Value = 472519.833804908  # Q11
$ws.Cells.Item(11, 18).Value = 7016939.429568958  # R11
$ws.Cells.Item(11, 19).Value = 10  # S11
$ws.Cells.Item(11, 20).Formula = '''Jämtland'  # T11
$ws.Cells.Item(11, 21).Formula = '''Krokom'  # U11
$ws.Cells.Item(11, 22).Formula = '''Jämtland'  # V11
$ws.Cells.Item(11, 23).Formula = '''Rödön'  # W11
$ws.Cells.Item(11, 24).Value = $null  # X11
$ws.Cells.Item(11, 25).Formula = '''2022-09-08'  # Y11
$ws.Cells.Item(11, 26).Formula = '''00:00'  # Z11
$ws.Cells.Item(11, 27).Formula = '''2022-09-08'  # AA11
$ws.Cells.Item(11, 28).Formula = '''00:00'  # AB11
$ws.Cells.Item(11, 29).Value = $null  # AC11
$ws.Cells.Item(11, 30).Value = $false  # AD11
$ws.Cells.Item(11, 31).Value = $false  # AE11
$ws.Cells.Item(11, 32).Value = $null  # AF11
$ws.Cells.Item(11, 33).Value = $false  # AG11
$ws.Cells.Item(11, 34).Value = $null  # AH11
$ws.Cells.Item(11, 35).Value = $null  # AI11
$ws.Cells.Item(11, 36).Value = $null  # AJ11
$ws.Cells.Item(11, 37).Value = $null  # AK11
$ws.Cells.Item(11, 38).Value = $null  # AL11
$ws.Cells.Item(11, 39).Value = $null  # AM11
$ws.Cells.Item(11, 40).Value = $null  # AN11
$ws.Cells.Item(11, 41).Value = $null  # AO11
$ws.Cells.Item(11, 42).Value = $null  # AP11
$ws.Cells.Item(11, 43).Value = $null  # AQ11
$ws.Cells.Item(11, 44).Value = $null  # AR11
$ws.Cells.Item(11, 45).Value = $null  # AS11
$ws.Cells.Item(11, 46).Formula = ''''  # AT11
$ws.Cells.Item(11, 47).Value = $null  # AU11
$ws.Cells.Item(11, 48).Value = $null  # AV11
$ws.Cells.Item(11, 49).Formula = '''Benny Öwre'  # AW11
$ws.Cells.Item(11, 50).Formula = '''Benny Öwre'  # AX11
$ws.Cells.Item(11, 51).Formula = ''''  # AY11

# ---- row 12 ----
$ws.Cells.Item(12, 1).Value = 106082248  # A12
$ws.Cells.Item(12, 2).Value = 56395  # B12
$ws.Cells.Item(12, 3).Formula = '''Ovaliderad'  # C12
$ws.Cells.Item(12, 4).Formula = '''NT'  # D12
$ws.Cells.Item(12, 5).Value = 100109  # E12
$ws.Cells.Item(12, 6).Formula = '''Tretåig hackspett'  # F12
$ws.Cells.Item(12, 7).Formula = '''Picoides tridactylus'  # G12
$ws.Cells.Item(12, 8).Formula = '''(Linnaeus, 1758)'  # H12
$ws.Cells.Item(12, 9).Formula = '''1'  # I12
$ws.Cells.Item(12, 10).Value = $null  # J12
$ws.Cells.Item(12, 11).Formula = ''''  # K12
$ws.Cells.Item(12, 12).Formula = ''''  # L12
$ws.Cells.Item(12, 13).Formula = '''födosökande'  # M12
$ws.Cells.Item(12, 14).Formula = ''''  # N12
$ws.Cells.Item(12, 15).Value = $null  # O12
$ws.Cells.Item(12, 16).Formula = '''Tretåig hackspett, Jmt'  # P12
$ws.Cells.Item(12, 17).Value = 472527.1055015869  # Q12
$ws.Cells.Item(12, 18).Value = 7016946.568404312  # R12
$ws.Cells.Item(12, 19).Value = 25  # S12
$ws.Cells.Item(12, 20).Formula = '''Jämtland'  # T12
$ws.Cells.Item(12, 21).Formula = '''Krokom'  # U12
$ws.Cells.Item(12, 22).Formula = '''Jämtland'  # V12
$ws.Cells.Item(12, 23).Formula = '''Rödön'  # W12
$ws.Cells.Item(12, 24).Value = $null  # X12
$ws.Cells.Item(12, 25).Formula = '''2023-01-20'  # Y12
$ws.Cells.Item(12, 26).Formula = '''11:00'  # Z12
$ws.Cells.Item(12, 27).Formula = '''2023-01-20'  # AA12
$ws.Cells.Item(12, 28).Formula = '''12:00'  # AB12
$ws.Cells.Item(12, 29).Formula = '''På död välbearbetad gran.'  # AC12
$ws.Cells.Item(12, 30).Value = $false  # AD12
$ws.Cells.Item(12, 31).Value = $false  # AE12
$ws.Cells.Item(12, 32).Value = $null  # AF12
$ws.Cells.Item(12, 33).Value = $false  # AG12
$ws.Cells.Item(12, 34).Value = $null  # AH12
$ws.Cells.Item(12, 35).Value = $null  # AI12
$ws.Cells.Item(12, 36).Value = $null  # AJ12
$ws.Cells.Item(12, 37).Value = $null  # AK12
$ws.Cells.Item(12, 38).Value = $null  # AL12
$ws.Cells.Item(12, 39).Value = $null  # AM12
$ws.Cells.Item(12, 40).Value = $null  # AN12
$ws.Cells.Item(12, 41).Value = $null  # AO12
$ws.Cells.Item(12, 42).Value = $null  # AP12
$ws.Cells.Item(12, 43).Value = $null  # AQ12
$ws.Cells.Item(12, 44).Value = $null  # AR12
$ws.Cells.Item(12, 45).Value = $null  # AS12
$ws.Cells.Item(12, 46).Formula = ''''  # AT12
$ws.Cells.Item(12, 47).Value = $null  # AU12
$ws.Cells.Item(12, 48).Value = $null  # AV12
$ws.Cells.Item(12, 49).Formula = '''Kristofer Holmsten'  # AW12
$ws.Cells.Item(12, 50).Formula = '''Kristofer Holmsten'  # AX12
$ws.Cells.Item(12, 51).Formula = ''''  # AY12

# ---- row 13 ----
$ws.Cells.Item(13, 1).Value = 103206713  # A13
$ws.Cells.Item(13, 2).Value = 96334  # B13
$ws.Cells.Item(13, 3).Formula = '''Ovaliderad'  # C13
$ws.Cells.Item(13, 4).Formula = '''VU'  # D13
$ws.Cells.Item(13, 5).Value = 220787  # E13
$ws.Cells.Item(13, 6).Formula = '''Knärot'  # F13
$ws.Cells.Item(13, 7).Formula = '''Goodyera repens'  # G13
$ws.Cells.Item(13, 8).Formula = '''(L.) R. Br.'  # H13
$ws.Cells.Item(13, 9).Formula = '''11'  # I13
$ws.Cells.Item(13, 10).Formula = '''stjälkar/strån/skott'  # J13
$ws.Cells.Item(13, 11).Formula = ''''  # K13
$ws.Cells.Item(13, 12).Formula = ''''  # L13
$ws.Cells.Item(13, 13).Value = $null  # M13
$ws.Cells.Item(13, 14).Formula = ''''  # N13
$ws.Cells.Item(13, 15).Value = $null  # O13
$ws.Cells.Item(13, 16).Formula = '''Bergom/Kroksgård - Rödön, Jmt'  # P13
$ws.Cells.Item(13, 17).Value = 472198.9007623708  # Q13
$ws.Cells.Item(13, 18).Value = 7017350.364024058  # R13
$ws.Cells.Item(13, 19).Value = 10  # S13
$ws.Cells.Item(13, 20).Formula = '''Jämtland'  # T13
$ws.Cells.Item(13, 21).Formula = '''Krokom'  # U13
$ws.Cells.Item(13, 22).Formula = '''Jämtland'  # V13
$ws.Cells.Item(13, 23).Formula = '''Rödön'  # W13
$ws.Cells.Item(13, 24).Value = $null  # X13
$ws.Cells.Item(13, 25).Formula = '''2022-08-22'  # Y13
$ws.Cells.Item(13, 26).Formula = '''00:00'  # Z13
$ws.Cells.Item(13, 27).Formula = '''2022-08-22'  # AA13
$ws.Cells.Item(13, 28).Formula = '''00:00'  # AB13
$ws.Cells.Item(13, 29).Value = $null  # AC13
$ws.Cells.Item(13, 30).Value = $false  # AD13
$ws.Cells.Item(13, 31).Value = $false  # AE13
$ws.Cells.Item(13, 32).Formula = ''''  # AF13
$ws.Cells.Item(13, 33).Value = $false  # AG13
$ws.Cells.Item(13, 34).Value = $null  # AH13
$ws.Cells.Item(13, 35).Value = $null  # AI13
$ws.Cells.Item(13, 36).Value = $null  # AJ13
$ws.Cells.Item(13, 37).Value = $null  # AK13
$ws.Cells.Item(13, 38).Value = $null  # AL13
$ws.Cells.Item(13, 39).Value = $null  # AM13
$ws.Cells.Item(13, 40).Value = $null  # AN13
$ws.Cells.Item(13, 41).Value = $null  # AO13
$ws.Cells.Item(13, 42).Value = $null  # AP13
$ws.Cells.Item(13, 43).Value = $null  # AQ13
$ws.Cells.Item(13, 44).Value = $null  # AR13
$ws.Cells.Item(13, 45).Value = $null  # AS13
$ws.Cells.Item(13, 46).Formula = ''''  # AT13
$ws.Cells.Item(13, 47).Value = $null  # AU13
$ws.Cells.Item(13, 48).Value = $null  # AV13
$ws.Cells.Item(13, 49).Formula = '''Benny Öwre'  # AW13
$ws.Cells.Item(13, 50).Formula = '''Benny Öwre'  # AX13
$ws.Cells.Item(13, 51).Formula = ''''  # AY13

# ---- row 14 ----
$ws.Cells.Item(14, 1).Value = 103636893  # A14
$ws.Cells.Item(14, 2).Value = 96334  # B14
$ws.Cells.Item(14, 3).Formula = '''Ovaliderad'  # C14
$ws.Cells.Item(14, 4).Formula = '''VU'  # D14
$ws.Cells.Item(14, 5).Value = 220787  # E14
$ws.Cells.Item(14, 6).Formula = '''Knärot'  # F14
$ws.Cells.Item(14, 7).Formula = '''Goodyera repens'  # G14
$ws.Cells.Item(14, 8).Formula = '''(L.) R. Br.'  # H14
$ws.Cells.Item(14, 9).Formula = ''''  # I14
$ws.Cells.Item(14, 10).Value = $null  # J14
$ws.Cells.Item(14, 11).Value = $null  # K14
$ws.Cells.Item(14, 12).Value = $null  # L14
$ws.Cells.Item(14, 13).Value = $null  # M14
$ws.Cells.Item(14, 14).Value = $null  # N14
$ws.Cells.Item(14, 15).Value = $null  # O14
$ws.Cells.Item(14, 16).Formula = '''Bergom - Rödön, Jmt'  # P14
$ws.Cells.Item(14, 17).Value = 472928.6702964447  # Q14
$ws.Cells.Item(14, 18).Value = 7016573.647136474  # R14
$ws.Cells.Item(14, 19).Value = 10  # S14
$ws.Cells.Item(14, 20).Formula = '''Jämtland'  # T14
$ws.Cells.Item(14, 21).Formula = '''Krokom'  # U14
$ws.Cells.Item(14, 22).Formula = '''Jämtland'  # V14
$ws.Cells.Item(14, 23).Formula = '''Rödön'  # W14
$ws.Cells.Item(14, 24).Value = $null  # X14
$ws.Cells.Item(14, 25).Formula = '''2022-09-08'  # Y14
$ws.Cells.Item(14, 26).Formula = '''00:00'  # Z14
$ws.Cells.Item(14, 27).Formula = '''2022-09-08'  # AA14
$ws.Cells.Item(14, 28).Formula = '''00:00'  # AB14
$ws.Cells.Item(14, 29).Value = $null  # AC14
$ws.Cells.Item(14, 30).Value = $false  # AD14
$ws.Cells.Item(14, 31).Value = $false  # AE14
$ws.Cells.Item(14, 32).Value = $null  # AF14
$ws.Cells.Item(14, 33).Value = $false  # AG14
$ws.Cells.Item(14, 34).Value = $null  # AH14
$ws.Cells.Item(14, 35).Value = $null  # AI14
$ws.Cells.Item(14, 36).Value = $null  # AJ14
$ws.Cells.Item(14, 37).Value = $null  # AK14
$ws.Cells.Item(14, 38).Value = $null  # AL14
$ws.Cells.Item(14, 39).Value = $null  # AM14
$ws.Cells.Item(14, 40).Value = $null  # AN14
$ws.Cells.Item(14, 41).Value = $null  # AO14
$ws.Cells.Item(14, 42).Value = $null  # AP14
$ws.Cells.Item(14, 43).Value = $null  # AQ14
$ws.Cells.Item(14, 44).Value = $null  # AR14
$ws.Cells.Item(14, 45).Value = $null  # AS14
$ws.Cells.Item(14, 46).Formula = ''''  # AT14
$ws.Cells.Item(14, 47).Value = $null  # AU14
$ws.Cells.Item(14, 48).Value = $null  # AV14
$ws.Cells.Item(14, 49).Formula = '''Benny Öwre'  # AW14
$ws.Cells.Item(14, 50).Formula = '''Benny Öwre'  # AX14
$ws.Cells.Item(14, 51).Formula = ''''  # AY14

# ---- row 15 ----
$ws.Cells.Item(15, 1).Value = 103636892  # A15
$ws.Cells.Item(15, 2).Value = 96334  # B15
$ws.Cells.Item(15, 3).Formula = '''Ovaliderad'  # C15
$ws.Cells.Item(15, 4).Formula = '''VU'  # D15
$ws.Cells.Item(15, 5).Value = 220787  # E15
$ws.Cells.Item(15, 6).Formula = '''Knärot'  # F15
$ws.Cells.Item(15, 7).Formula = '''Goodyera repens'  # G15
$ws.Cells.Item(15, 8).Formula = '''(L.) R. Br.'  # H15
$ws.Cells.Item(15, 9).Formula = ''''  # I15
$ws.Cells.Item(15, 10).Value = $null  # J15
$ws.Cells.Item(15, 11).Value = $null  # K15
$ws.Cells.Item(15, 12).Value = $null  # L15
$ws.Cells.Item(15, 13).Value = $null  # M15
$ws.Cells.Item(15, 14).Value = $null  # N15
$ws.Cells.Item(15, 15).Value = $null  # O15
$ws.Cells.Item(15, 16).Formula = '''Bergom - Rödön, Jmt'  # P15
$ws.Cells.Item(15, 17).Value = 472939.4717169611  # Q15
$ws.Cells.Item(15, 18).Value = 7016571.755861398  # R15
$ws.Cells.Item(15, 19).Value = 10  # S15
$ws.Cells.Item(15, 20).Formula = '''Jämtland'  # T15
$ws.Cells.Item(15, 21).Formula = '''Krokom'  # U15
$ws.Cells.Item(15, 22).Formula = '''Jämtland'  # V15
$ws.Cells.Item(15, 23).Formula = '''Rödön'  # W15
$ws.Cells.Item(15, 24).Value = $null  # X15
$ws.Cells.Item(15, 25).Formula = '''2022-09-08'  # Y15
$ws.Cells.Item(15, 26).Formula = '''00:00'  # Z15
$ws.Cells.Item(15, 27).Formula = '''2022-09-08'  # AA15
$ws.Cells.Item(15, 28).Formula = '''00:00'  # AB15
$ws.Cells.Item(15, 29).Value = $null  # AC15
$ws.Cells.Item(15, 30).Value = $false  # AD15
$ws.Cells.Item(15, 31).Value = $false  # AE15
$ws.Cells.Item(15, 32).Value = $null  # AF15
$ws.Cells.Item(15, 33).Value = $false  # AG15
$ws.Cells.Item(15, 34).Value = $null  # AH15
$ws.Cells.Item(15, 35).Value = $null  # AI15
$ws.Cells.Item(15, 36).Value = $null  # AJ15
$ws.Cells.Item(15, 37).Value = $null  # AK15
$ws.Cells.Item(15, 38).Value = $null  # AL15
$ws.Cells.Item(15, 39).Value = $null  # AM15
$ws.Cells.Item(15, 40).Value = $null  # AN15
$ws.Cells.Item(15, 41).Value = $null  # AO15
$ws.Cells.Item(15, 42).Value = $null  # AP15
$ws.Cells.Item(15, 43).Value = $null  # AQ15
$ws.Cells.Item(15, 44).Value = $null  # AR15
$ws.Cells.Item(15, 45).Value = $null  # AS15
$ws.Cells.Item(15, 46).Formula = ''''  # AT15
$ws.Cells.Item(15, 47).Value = $null  # AU15
$ws.Cells.Item(15, 48).Value = $null  # AV15
$ws.Cells.Item(15, 49).Formula = '''Benny Öwre'  # AW15
$ws.Cells.Item(15, 50).Formula = '''Benny Öwre'  # AX15
$ws.Cells.Item(15, 51).Formula = ''''  # AY15

